# Error Calculations and Plots
# Apply the edits described by the diff:
#  1. Remove two rows that were excluded from the dataset (old "RM 232" at
#     row 26, and old "SC 92" at what becomes row 27 after the first delete).
#     Excel shifts everything below up automatically, just like a manual
#     right-click > Delete on the row headers.
#  2. Update individual "D" column (imputed) values that changed for the
#     remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two removed rows -----------------------------------
$ws.Rows.Item(26).Delete()   # removes old row 26 ("RM 232"); rows below shift up by 1
$ws.Rows.Item(27).Delete()   # old row 28 ("SC 92") is now at row 27; remove it too

# --- 2. Fix up individual cell values (rows 1-25 keep their row numbers) --
$ws.Range("E2").Value = -7.2
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = 17.97
$ws.Range("F5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = ""
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E22").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = 16.48

# --- 3. Fix up individual cell values in the rows that shifted up --------
# (row numbers below refer to the NEW row positions, after the deletions)
$ws.Range("F27").Value = ""        # was row 29 "SC 101": F 17 -> blank
$ws.Range("F29").Value = 18.06     # was row 31 "SC 119": F blank -> 18.06
$ws.Range("C30").Value = 11.4      # was row 32 "SC 120": C blank -> 11.4
$ws.Range("E31").Value = -8.1      # was row 33 "SC 132": E blank -> -8.1
$ws.Range("C32").Value = ""        # was row 34 "SC 193": C 10.5 -> blank
$ws.Range("E33").Value = -10.7     # was row 35 "SC 232": E blank -> -10.7
